$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Recorded By" (column G) email lists reordered across many rows.
#    Content only - no formatting changes.
# ---------------------------------------------------------------------------

$group1 = "shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
foreach ($r in @(2, 21, 40)) {
    $ws.Range("G$r").Value = $group1
}

$ws.Range("G9").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

$group3 = "wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
foreach ($r in @(18, 37, 56, 75, 94, 113)) {
    $ws.Range("G$r").Value = $group3
}

$group4 = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
foreach ($r in @(19, 76, 95)) {
    $ws.Range("G$r").Value = $group4
}

$group5 = "neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
foreach ($r in @(20, 38, 39, 57, 58, 77, 96, 115)) {
    $ws.Range("G$r").Value = $group5
}

$group6 = "AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
foreach ($r in @(24, 81)) {
    $ws.Range("G$r").Value = $group6
}

$group7 = "norhan.mohamed@med.asu.edu.eg, yasmintarek@med.asu.edu.eg"
foreach ($r in @(25, 82)) {
    $ws.Range("G$r").Value = $group7
}

$ws.Range("G28").Value = "arwaelsayed03@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"

$group9 = "Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
foreach ($r in @(43, 100)) {
    $ws.Range("G$r").Value = $group9
}

$group10 = "esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
foreach ($r in @(47, 104)) {
    $ws.Range("G$r").Value = $group10
}

$group11 = "servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
foreach ($r in @(59, 78, 97)) {
    $ws.Range("G$r").Value = $group11
}

$group12 = "yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
foreach ($r in @(66, 85)) {
    $ws.Range("G$r").Value = $group12
}

# ---------------------------------------------------------------------------
# 2) Summary statistics that shift as a consequence of the attendance edits.
# ---------------------------------------------------------------------------

$ws.Range("L6").Value = 41
$ws.Range("L8").Value = 70
$ws.Range("L9").Value = "36.0%"
$ws.Range("L10").Value = "45.9%"

$ws.Range("O15").Value = 8
$ws.Range("P15").Value = 0
$ws.Range("R15").Value = "42.1%"
$ws.Range("S15").Value = "41.8%"

$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 11

# ---------------------------------------------------------------------------
# 3) Row 35 ("PHARMACOLOGY" A2 session on 30/11/2025) flips from Pending to
#    Not Recorded - copy the "Not Recorded" look (pink highlight) from a
#    row that already carries it (row 62), then update the status text.
# ---------------------------------------------------------------------------

$ws.Range("A62:I62").Copy() | Out-Null
$ws.Range("A35:I35").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I35").Value = "Not Recorded"

# ---------------------------------------------------------------------------
# 4) Row 16 ("PHARMACOLOGY" A1 session on 30/11/2025) flips from Not
#    Recorded to Recorded - copy the "Recorded" look (green highlight) from
#    a row that already carries it (row 2), then update content.
# ---------------------------------------------------------------------------

$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A16:I16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("G16").Value = "nourhan.mostafa@med.asu.edu.eg"
$ws.Range("H16").Value = "7/203"
$ws.Range("I16").Value = "Recorded"
